# ---------------------------------------------------------------------------
# Jason Puzon Dy (Q0501) workbook update
#   - "adding averages and more checks"
#   - Training Dashboard: refresh PERIOD TO EXPIRE / LAST UPDATE columns
#   - Exam Dashboard: update comments + shrink the COMMENTS column width
#   - Shared style table: headers/title now use a bold white-on-blue look
# ---------------------------------------------------------------------------

function Set-LiteralText {
    # Writes $text into $range as a plain literal value (never auto-parsed
    # into a date/number by the smart-entry logic), while leaving the
    # range's existing number format / style completely untouched.
    param($range, [string]$text)

    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $value = $range.Value
    $range.Value = $value
}

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: PERIOD TO EXPIRE (H) / LAST UPDATE (I) -----------

$ws1.Range("H3").Value = 244
Set-LiteralText $ws1.Range("I3") "16-Sep-2025"

$ws1.Range("H4").Value = 219
Set-LiteralText $ws1.Range("I4") "16-Sep-2025"

$ws1.Range("H5").Value = 255
Set-LiteralText $ws1.Range("I5") "16-Sep-2025"

$ws1.Range("H6").Value = 254
Set-LiteralText $ws1.Range("I6") "16-Sep-2025"

$ws1.Range("H7").Value = -103
Set-LiteralText $ws1.Range("I7") "16-Sep-2025"

$ws1.Range("H8").Value = -343
Set-LiteralText $ws1.Range("I8") "16-Sep-2025"

$ws1.Range("H9").Value = -50
Set-LiteralText $ws1.Range("I9") "16-Sep-2025"

# --- Exam Dashboard: comments + column width -------------------------------

Set-LiteralText $ws2.Range("E3") "date is valid"
Set-LiteralText $ws2.Range("E4") "date is valid"

# Column E (COMMENTS) narrows from 44 to 15 characters; copy the exact
# COM-reported width from a column that is already 15 wide so the saved
# <col width="..."> lines up exactly instead of drifting from rounding.
$ws2.Columns.Item(5).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth

# --- Shared look for the report titles and table headers ------------------
# Title (merged A1) and header row both become bold white text; the title
# also drops its old 14pt override so it matches the header's font.

foreach ($pair in @(
        @{ Sheet = $ws1; Title = "A1"; Header = "A2:K2" },
        @{ Sheet = $ws2; Title = "A1"; Header = "A2:G2" }
    )) {
    $sheet = $pair.Sheet

    $titleFont = $sheet.Range($pair.Title).Font
    $titleFont.Bold = $true
    $titleFont.Size = 11
    $titleFont.Color = 0xFFFFFF

    $headerFont = $sheet.Range($pair.Header).Font
    $headerFont.Bold = $true
    $headerFont.Color = 0xFFFFFF
}
